$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing column B ("Jun_13" data),
# shifting the old B column to D and the old C column to E.
$ws.Range("B:C").Insert()

# New column C (Jun_15) and column B (Jun_17) headers -- set C before B
# so the new shared strings are appended in "Jun_15", "Jun_17" order.
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the two new columns (rows 2-27) with the ticker value, copied
# from the original ticker column (now shifted to column D).
for ($r = 2; $r -le 27; $r++) {
    $ticker = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 2).Value = $ticker
    $ws.Cells.Item($r, 3).Value = $ticker
}

# Match the original column's custom width (8 characters) on the new
# columns and the shifted-over original column.
$ws.Columns("C").ColumnWidth = 8.0
$ws.Columns("D").ColumnWidth = 8.0
$ws.Columns("E").ColumnWidth = 8.0
